$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 72.31432066666666
$ws.Range("H2").Value = 216.942962
$ws.Range("I2").Value = 0.1293883843050027
$ws.Range("J2").Value = 0.1293883843050027
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 123.0738876674464
$ws.Range("R2").Value = 1107.664989007018
$ws.Range("S2").Value = 0.002724431433220291
$ws.Range("T2").Value = 0.002724431433220291

$ws.Range("G3").Value = 72.31432066666666
$ws.Range("H3").Value = 216.942962
$ws.Range("I3").Value = 0.1293883843050027
$ws.Range("J3").Value = 0.1293883843050027
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("Q3").Value = 4519.930466290338
$ws.Range("R3").Value = 40679.37419661304
$ws.Range("S3").Value = 0.1000556728296852
$ws.Range("T3").Value = 0.1000556728296852

$ws.Range("G4").Value = 72.31432066666666
$ws.Range("H4").Value = 216.942962
$ws.Range("I4").Value = 0.1293883843050027
$ws.Range("J4").Value = 0.1293883843050027
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 31.374941993326
$ws.Range("R4").Value = 282.374477939934
$ws.Range("S4").Value = 0.0006945330142901642
$ws.Range("T4").Value = 0.0006945330142901643

$ws.Range("G5").Value = 72.31432066666666
$ws.Range("H5").Value = 216.942962
$ws.Range("I5").Value = 0.1293883843050027
$ws.Range("J5").Value = 0.1293883843050027
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 1155.032098387136
$ws.Range("R5").Value = 10395.28888548422
$ws.Range("S5").Value = 0.0255684273477017
$ws.Range("T5").Value = 0.02556842734770171

$ws.Range("G6").Value = 72.31432066666666
$ws.Range("H6").Value = 216.942962
$ws.Range("I6").Value = 0.1293883843050027
$ws.Range("J6").Value = 0.1293883843050027
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 15.59952473034556
$ws.Range("R6").Value = 140.39572257311
$ws.Range("S6").Value = 0.0003453196801054014
$ws.Range("T6").Value = 0.0003453196801054014

$ws.Range("I7").Value = 0.2756445674916043
$ws.Range("J7").Value = 0.2756445674916043
$ws.Range("M7").Value = 1.701929666666667
$ws.Range("N7").Value = 5.105789
$ws.Range("O7").Value = 0.02105622887134972
$ws.Range("P7").Value = 0.02105622887134972
$ws.Range("Q7").Value = 262.1923808526288
$ws.Range("R7").Value = 2359.731427673659
$ws.Range("S7").Value = 0.005804035100247425
$ws.Range("T7").Value = 0.005804035100247425

$ws.Range("I8").Value = 0.2756445674916043
$ws.Range("J8").Value = 0.2756445674916043
$ws.Range("O8").Value = 0.7732971809418951
$ws.Range("P8").Value = 0.7732971809418953
$ws.Range("S8").Value = 0.2131551669832056
$ws.Range("T8").Value = 0.2131551669832056

$ws.Range("I9").Value = 0.2756445674916043
$ws.Range("J9").Value = 0.2756445674916043
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4338690000000001
$ws.Range("N9").Value = 1.301607
$ws.Range("O9").Value = 0.005367815805265532
$ws.Range("P9").Value = 0.005367815805265533
$ws.Range("Q9").Value = 66.84009822271301
$ws.Range("R9").Value = 601.5608840044172
$ws.Range("S9").Value = 0.001479609266017015
$ws.Range("T9").Value = 0.001479609266017015

$ws.Range("I10").Value = 0.2756445674916043
$ws.Range("J10").Value = 0.2756445674916043
$ws.Range("M10").Value = 15.972384
$ws.Range("N10").Value = 47.917152
$ws.Range("O10").Value = 0.1976099128607259
$ws.Range("P10").Value = 0.1976099128607259
$ws.Range("Q10").Value = 2460.640689726368
$ws.Range("R10").Value = 22145.76620753732
$ws.Range("S10").Value = 0.05447009896254841
$ws.Range("T10").Value = 0.05447009896254842

$ws.Range("I11").Value = 0.2756445674916043
$ws.Range("J11").Value = 0.2756445674916043
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2157183333333333
$ws.Range("N11").Value = 0.647155
$ws.Range("O11").Value = 0.002668861520763652
$ws.Range("P11").Value = 0.002668861520763652
$ws.Range("Q11").Value = 33.23269140786723
$ws.Range("R11").Value = 299.094222670805
$ws.Range("S11").Value = 0.0007356571795858823
$ws.Range("T11").Value = 0.0007356571795858823

$ws.Range("G12").Value = 168.0546723333333
$ws.Range("H12").Value = 504.1640170000001
$ws.Range("I12").Value = 0.300691790058393
$ws.Range("J12").Value = 0.300691790058393
$ws.Range("M12").Value = 1.701929666666667
$ws.Range("N12").Value = 5.105789
$ws.Range("O12").Value = 0.02105622887134972
$ws.Range("P12").Value = 0.02105622887134972
$ws.Range("Q12").Value = 286.0172324660459
$ws.Range("R12").Value = 2574.155092194413
$ws.Range("S12").Value = 0.006331435151205362
$ws.Range("T12").Value = 0.006331435151205362

$ws.Range("G13").Value = 168.0546723333333
$ws.Range("H13").Value = 504.1640170000001
$ws.Range("I13").Value = 0.300691790058393
$ws.Range("J13").Value = 0.300691790058393
$ws.Range("O13").Value = 0.7732971809418951
$ws.Range("P13").Value = 0.7732971809418953
$ws.Range("Q13").Value = 10504.08033262503
$ws.Range("R13").Value = 94536.72299362531
$ws.Range("S13").Value = 0.2325241135845275
$ws.Range("T13").Value = 0.2325241135845275

$ws.Range("G14").Value = 168.0546723333333
$ws.Range("H14").Value = 504.1640170000001
$ws.Range("I14").Value = 0.300691790058393
$ws.Range("J14").Value = 0.300691790058393
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4338690000000001
$ws.Range("N14").Value = 1.301607
$ws.Range("O14").Value = 0.005367815805265532
$ws.Range("P14").Value = 0.005367815805265533
$ws.Range("Q14").Value = 72.91371263059102
$ws.Range("R14").Value = 656.2234136753192
$ws.Range("S14").Value = 0.001614058143189027
$ws.Range("T14").Value = 0.001614058143189027

$ws.Range("G15").Value = 168.0546723333333
$ws.Range("H15").Value = 504.1640170000001
$ws.Range("I15").Value = 0.300691790058393
$ws.Range("J15").Value = 0.300691790058393
$ws.Range("M15").Value = 15.972384
$ws.Range("N15").Value = 47.917152
$ws.Range("O15").Value = 0.1976099128607259
$ws.Range("P15").Value = 0.1976099128607259
$ws.Range("Q15").Value = 2684.233759502176
$ws.Range("R15").Value = 24158.10383551959
$ws.Range("S15").Value = 0.05941967843137472
$ws.Range("T15").Value = 0.05941967843137473

$ws.Range("G16").Value = 168.0546723333333
$ws.Range("H16").Value = 504.1640170000001
$ws.Range("I16").Value = 0.300691790058393
$ws.Range("J16").Value = 0.300691790058393
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2157183333333333
$ws.Range("N16").Value = 0.647155
$ws.Range("O16").Value = 0.002668861520763652
$ws.Range("P16").Value = 0.002668861520763652
$ws.Range("Q16").Value = 36.25247382462612
$ws.Range("R16").Value = 326.272264421635
$ws.Range("S16").Value = 0.0008025047480963876
$ws.Range("T16").Value = 0.0008025047480963876

$ws.Range("G17").Value = 65.818911
$ws.Range("H17").Value = 197.456733
$ws.Range("I17").Value = 0.11776647381174
$ws.Range("J17").Value = 0.11776647381174
$ws.Range("M17").Value = 1.701929666666667
$ws.Range("N17").Value = 5.105789
$ws.Range("O17").Value = 0.02105622887134972
$ws.Range("P17").Value = 0.02105622887134972
$ws.Range("Q17").Value = 112.019157258593
$ws.Range("R17").Value = 1008.172415327337
$ws.Range("S17").Value = 0.002479717825951811
$ws.Range("T17").Value = 0.002479717825951811

$ws.Range("G18").Value = 65.818911
$ws.Range("H18").Value = 197.456733
$ws.Range("I18").Value = 0.11776647381174
$ws.Range("J18").Value = 0.11776647381174
$ws.Range("O18").Value = 0.7732971809418951
$ws.Range("P18").Value = 0.7732971809418953
$ws.Range("Q18").Value = 4113.941724741717
$ws.Range("R18").Value = 37025.47552267546
$ws.Range("S18").Value = 0.09106848220808608
$ws.Range("T18").Value = 0.0910684822080861

$ws.Range("G19").Value = 65.818911
$ws.Range("H19").Value = 197.456733
$ws.Range("I19").Value = 0.11776647381174
$ws.Range("J19").Value = 0.11776647381174
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.4338690000000001
$ws.Range("N19").Value = 1.301607
$ws.Range("O19").Value = 0.005367815805265532
$ws.Range("P19").Value = 0.005367815805265533
$ws.Range("Q19").Value = 28.556785096659
$ws.Range("R19").Value = 257.011065869931
$ws.Range("S19").Value = 0.0006321487394570475
$ws.Range("T19").Value = 0.0006321487394570476

$ws.Range("G20").Value = 65.818911
$ws.Range("H20").Value = 197.456733
$ws.Range("I20").Value = 0.11776647381174
$ws.Range("J20").Value = 0.11776647381174
$ws.Range("M20").Value = 15.972384
$ws.Range("N20").Value = 47.917152
$ws.Range("O20").Value = 0.1976099128607259
$ws.Range("P20").Value = 0.1976099128607259
$ws.Range("Q20").Value = 1051.284920953824
$ws.Range("R20").Value = 9461.564288584417
$ws.Range("S20").Value = 0.02327182262785291
$ws.Range("T20").Value = 0.02327182262785291

$ws.Range("G21").Value = 65.818911
$ws.Range("H21").Value = 197.456733
$ws.Range("I21").Value = 0.11776647381174
$ws.Range("J21").Value = 0.11776647381174
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2157183333333333
$ws.Range("N21").Value = 0.647155
$ws.Range("O21").Value = 0.002668861520763652
$ws.Range("P21").Value = 0.002668861520763652
$ws.Range("Q21").Value = 14.198345782735
$ws.Range("R21").Value = 127.785112044615
$ws.Range("S21").Value = 0.0003143024103921733
$ws.Range("T21").Value = 0.0003143024103921733

$ws.Range("G22").Value = 98.64960366666666
$ws.Range("H22").Value = 295.948811
$ws.Range("I22").Value = 0.17650878433326
$ws.Range("J22").Value = 0.1765087843332599
$ws.Range("M22").Value = 1.701929666666667
$ws.Range("N22").Value = 5.105789
$ws.Range("O22").Value = 0.02105622887134972
$ws.Range("P22").Value = 0.02105622887134972
$ws.Range("Q22").Value = 167.8946870852088
$ws.Range("R22").Value = 1511.052183766879
$ws.Range("S22").Value = 0.003716609360724829
$ws.Range("T22").Value = 0.003716609360724829

$ws.Range("G23").Value = 98.64960366666666
$ws.Range("H23").Value = 295.948811
$ws.Range("I23").Value = 0.17650878433326
$ws.Range("J23").Value = 0.1765087843332599
$ws.Range("O23").Value = 0.7732971809418951
$ws.Range("P23").Value = 0.7732971809418953
$ws.Range("Q23").Value = 6165.989599152339
$ws.Range("R23").Value = 55493.90639237105
$ws.Range("S23").Value = 0.1364937453363909
$ws.Range("T23").Value = 0.1364937453363909

$ws.Range("G24").Value = 98.64960366666666
$ws.Range("H24").Value = 295.948811
$ws.Range("I24").Value = 0.17650878433326
$ws.Range("J24").Value = 0.1765087843332599
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 0.4338690000000001
$ws.Range("N24").Value = 1.301607
$ws.Range("O24").Value = 0.005367815805265532
$ws.Range("P24").Value = 0.005367815805265533
$ws.Range("Q24").Value = 42.80100489325301
$ws.Range("R24").Value = 385.209044039277
$ws.Range("S24").Value = 0.000947466642312278
$ws.Range("T24").Value = 0.000947466642312278

$ws.Range("G25").Value = 98.64960366666666
$ws.Range("H25").Value = 295.948811
$ws.Range("I25").Value = 0.17650878433326
$ws.Range("J25").Value = 0.1765087843332599
$ws.Range("M25").Value = 15.972384
$ws.Range("N25").Value = 47.917152
$ws.Range("O25").Value = 0.1976099128607259
$ws.Range("P25").Value = 0.1976099128607259
$ws.Range("Q25").Value = 1575.669351211808
$ws.Range("R25").Value = 14181.02416090627
$ws.Range("S25").Value = 0.03487988549124817
$ws.Range("T25").Value = 0.03487988549124817

$ws.Range("G26").Value = 98.64960366666666
$ws.Range("H26").Value = 295.948811
$ws.Range("I26").Value = 0.17650878433326
$ws.Range("J26").Value = 0.1765087843332599
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.2157183333333333
$ws.Range("N26").Value = 0.647155
$ws.Range("O26").Value = 0.002668861520763652
$ws.Range("P26").Value = 0.002668861520763652
$ws.Range("Q26").Value = 21.28052808696722
$ws.Range("R26").Value = 191.524752782705
$ws.Range("S26").Value = 0.0004710775025838078
$ws.Range("T26").Value = 0.0004710775025838077
